# Apply updated cryptocurrency price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Route through a formula + paste-as-values round trip so that
    # numeric-looking strings (e.g. "1.00", "7.10") are stored as
    # literal text instead of being coerced into numbers, while
    # leaving the cells existing style/format untouched.
    $escaped = $text -replace '"', '""'
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

$ws.Range('D2').Value = '67.754.03'
$ws.Range('E2').Value = '  +1.24%  '

$ws.Range('D3').Value = '3.538.46'
$ws.Range('E3').Value = '  +0.17%  '

Set-TextValue $ws.Range('D4') '1.00'
$ws.Range('E4').Value = '  -0.04%  '

Set-TextValue $ws.Range('D5') '610.85'
$ws.Range('E5').Value = '  +0.66%  '

Set-TextValue $ws.Range('D6') '152.02'
$ws.Range('E6').Value = '  -1.31%  '

$ws.Range('D7').Value = '3.538.00'
$ws.Range('E7').Value = '  +0.22%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  -0.63%  '

Set-TextValue $ws.Range('D10') '0.139'
$ws.Range('E10').Value = '  -1.34%  '

Set-TextValue $ws.Range('D11') '7.10'
$ws.Range('E11').Value = '  +4.02%  '

$ws.Range('E12').Value = '  -1.02%  '

Set-TextValue $ws.Range('D13') '0.0000218'
$ws.Range('E13').Value = '  -1.52%  '

$ws.Range('D14').Value = '4.139.76'
$ws.Range('E14').Value = '  +0.25%  '

Set-TextValue $ws.Range('D15') '31.92'
$ws.Range('E15').Value = '  -0.23%  '

$ws.Range('D16').Value = '3.552.37'
$ws.Range('E16').Value = '  +0.98%  '

$ws.Range('D17').Value = '67.566.18'
$ws.Range('E17').Value = '  +0.92%  '

$ws.Range('E18').Value = '  -0.76%  '

Set-TextValue $ws.Range('D19') '6.38'
$ws.Range('E19').Value = '  +0.27%  '

Set-TextValue $ws.Range('D20') '15.20'
$ws.Range('E20').Value = '  -1.62%  '

Set-TextValue $ws.Range('D21') '9.71'
$ws.Range('E21').Value = '  +4.09%  '

Set-TextValue $ws.Range('D22') '445.81'
$ws.Range('E22').Value = '  -0.93%  '

Set-TextValue $ws.Range('D23') '0.621'
$ws.Range('E23').Value = '  -2.58%  '

Set-TextValue $ws.Range('D24') '76.99'
$ws.Range('E24').Value = '  -2.58%  '

Set-TextValue $ws.Range('D25') '0.0000129'
$ws.Range('E25').Value = '  +4.71%  '

$ws.Range('D26').Value = '3.682.84'
$ws.Range('E26').Value = '  +0.23%  '

$ws.Range('E27').Value = '  +0.11%  '

Set-TextValue $ws.Range('D28') '10.17'
$ws.Range('E28').Value = '  -0.88%  '

Set-TextValue $ws.Range('D29') '8.64'
$ws.Range('E29').Value = '  +3.89%  '

$ws.Range('E30').Value = '  -0.40%  '

$ws.Range('E31').Value = '  -3.52%  '

$ws.Range('E32').Value = '  +6.97%  '

Set-TextValue $ws.Range('D34') '25.66'
$ws.Range('E34').Value = '  -0.95%  '

Set-TextValue $ws.Range('D35') '6.19'
$ws.Range('E35').Value = '  +0.14%  '

$ws.Range('D36').Value = '3.527.86'
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('E37').Value = '  -2.59%  '

Set-TextValue $ws.Range('D38') '8.03'
$ws.Range('E38').Value = '  -0.98%  '

$ws.Range('E40').Value = '  +0.01%  '

Set-TextValue $ws.Range('D41') '175.94'
$ws.Range('E41').Value = '  -0.08%  '

Set-TextValue $ws.Range('D42') '2.18'
$ws.Range('E42').Value = '  +1.85%  '

Set-TextValue $ws.Range('D43') '0.0893'
$ws.Range('E43').Value = '  +2.17%  '

$ws.Range('E44').Value = '  -3.84%  '

Set-TextValue $ws.Range('D45') '0.887'
$ws.Range('E45').Value = '  -0.45%  '

Set-TextValue $ws.Range('D46') '28.58'
$ws.Range('E46').Value = '  +1.25%  '

Set-TextValue $ws.Range('D47') '45.49'

Set-TextValue $ws.Range('D48') '2.66'
$ws.Range('E48').Value = '  -0.61%  '

$ws.Range('E49').Value = '  +3.34%  '

Set-TextValue $ws.Range('D50') '7.59'
$ws.Range('E50').Value = '  -0.77%  '

$ws.Range('E51').Value = '  +0.27%  '
